$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Loans")

# Row 6: new loan record (mirrors the shape of row 5 - a returned loan)
$ws.Range("A6").Value = "'0"
$ws.Range("B6").Value = "'1"
$ws.Range("C6").Value = "Goku"
$ws.Range("D6").Value = "Wed May 23 02:06:20 ART 2018"
$ws.Range("E6").Value = "Thu May 24 02:06:20 ART 2018"

# Drop the auto-applied quote-prefix formatting so A6/B6 stay on the
# workbook's default style, matching the plain text-number cells used
# elsewhere in this sheet (e.g. A5/B5).
$ws.Range("A6:B6").Style = "Normal"
